$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of J2:J11 (k column)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# New summary rows 14-17
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Apply bold/12pt/vertical-center style to the new B column summary cells
$summaryRange = $ws.Range("B14:B17")
$summaryFont = $summaryRange.Font
$summaryFont.Bold = $true
$summaryFont.Size = 12
$summaryRange.VerticalAlignment = -4108

# Selection as referenced in the diff
$ws.Range("J2:J12").Select()
